$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "27.978.45"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "  -0.60%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.764.89"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "  -2.83%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.003"
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "  +0.50%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "338.79"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "  +0.12%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.9972"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "  +0.11%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.3761"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "  -4.50%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3365"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "  -3.69%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "45.97"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "  -5.10%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "1.130"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "  -5.93%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07215"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "  -5.08%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "22.75"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "  +2.49%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.9983"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "  +0.06%  "

$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "  -5.13%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "7.191"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "  -0.24%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "1.761.69"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "  -2.88%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.00001057"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "  -4.58%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.06579"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "  -1.97%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "80.65"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "  -5.62%  "

$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "  +0.17%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "16.99"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "  -4.88%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.291"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "  -4.33%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "28.022.27"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "  -0.39%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "11.75"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "  -8.58%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.375"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "  -1.49%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "153.28"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "  -1.04%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "2.350"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "  -8.83%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "19.79"
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "  -7.59%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.282"
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = "  -15.83%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.962.90"
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = "  -2.81%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "130.83"
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = "  -3.82%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "4.018"
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = "  -0.53%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "5.850"
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = "  -5.07%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.08769"
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = "  -1.05%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "12.28"
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = "  -7.81%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.02348"
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = "  -3.44%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.6618"
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = "  -4.86%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.06235"
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "  -4.99%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "5.158"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "  -6.79%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.2116"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "  -5.30%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "1.216"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "  -4.09%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "1.451"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "  -10.24%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "8.064"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "  -5.76%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.9972"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "  +0.12%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "13.71"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "  -6.36%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "3.833"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "  -1.10%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.6051"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "  -7.14%  "

$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "  -1.84%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "2.020"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "  -7.09%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.07243"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "  +0.25%  "

$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "  +1.33%  "

